$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '301.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.50%'
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '6'
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '8.03%'
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '6'
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.012'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.23%'
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '6'
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07856'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.38%'
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '6'
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.212'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-6.13%'
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '6'
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.006'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.01%'
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '6'
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.010'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.55%'
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '6'
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9102'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.93%'
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '6'
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1856'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.32%'
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '6'
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09219'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-8.68%'
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '6'
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08488'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.64%'
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '6'
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03519'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '5.90%'
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '6'
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09958'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.64%'
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '6'
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001470'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-2.16%'
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '6'
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005652'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.77%'
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '6'
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.471'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.18%'
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '6'
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.061'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.63%'
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '6'
$ws.Range("G18").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.87%'
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '6'
$ws.Range("G19").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.03%'
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '6'
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.552'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '5.42%'
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '6'
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2229'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-6.73%'
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '6'
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04645'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.67%'
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '6'
$ws.Range("G23").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.88%'
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '6'
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004449'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.09%'
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '6'
$ws.Range("G25").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.20%'
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '6'
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004747'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '39.90%'
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '6'
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '6'
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '6'
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '6'
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '6'
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '6'
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '6'
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '6'
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '6'
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '6'
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '6'
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '6'
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01754'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-1.38%'
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '6'
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04716'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.27%'
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '6'
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007909'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.12%'
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '6'
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1390'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.33%'
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '6'
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007666'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '10.30%'
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '6'
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002188'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '2.13%'
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '6'
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01037'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '12.74%'
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '6'
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00005992'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-2.19%'
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '6'
$ws.Range("G46").Style = "Normal"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '6'
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.672'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '235.86%'
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '6'
$ws.Range("G48").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '34.81%'
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '6'
$ws.Range("G49").Style = "Normal"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '6'
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001999'
$ws.Range("D51").Style = "Normal"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '6'
$ws.Range("G51").Style = "Normal"
